$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1272.5
$ws.Range("I8").Value = 30
$ws.Range("J8").Value = 5000
$ws.Range("K8").Value = 90
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 49
$ws.Range("N8").Value = -15278

$ws.Range("H41").Value = 166.16667
$ws.Range("I41").Value = 87
$ws.Range("J41").Value = 182
$ws.Range("K41").Value = 87
$ws.Range("L41").Value = 182
$ws.Range("M41").Value = 353
$ws.Range("N41").Value = -1062

$ws.Range("H55").Value = 163.63637
$ws.Range("I55").Value = 124.85714
$ws.Range("J55").Value = 231.5
$ws.Range("K55").Value = 124.85714
$ws.Range("L55").Value = 231.5
$ws.Range("M55").Value = 89.14286
$ws.Range("N55").Value = -659.5

$ws.Range("H80").Value = 470.58334
$ws.Range("I80").Value = 274.93332
$ws.Range("J80").Value = 796.6667
$ws.Range("K80").Value = 824.7999599999999
$ws.Range("L80").Value = 2390.0001
$ws.Range("M80").Value = 173.2000400000001
$ws.Range("N80").Value = -4386.0001

$ws.Range("H83").Value = 470.58334
$ws.Range("I83").Value = 274.93332
$ws.Range("J83").Value = 796.6667
$ws.Range("K83").Value = 2474.39988
$ws.Range("L83").Value = 7170.0003
$ws.Range("M83").Value = 2517.60012
$ws.Range("N83").Value = -17154.0003

$ws.Range("H129").Value = 1341.5082
$ws.Range("J129").Value = 1398.8928
$ws.Range("L129").Value = 4196.678400000001
$ws.Range("N129").Value = -14196.6784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5225.7144
$ws.Range("I2").Value = 1750
$ws.Range("K2").Value = 1750
$ws.Range("M2").Value = -1637

$ws.Range("H116").Value = 5225.7144
$ws.Range("I116").Value = 1750
$ws.Range("K116").Value = 1750
$ws.Range("M116").Value = 544

$ws.Range("H122").Value = 1969
$ws.Range("I122").Value = 1900.6875
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 5702.0625
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = -3252.0625
$ws.Range("N122").Value = -11899.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5225.7144
$ws.Range("I3").Value = 1750
$ws.Range("K3").Value = 1750
$ws.Range("M3").Value = -1636

$ws.Range("H86").Value = 1835.4546
$ws.Range("I86").Value = 1990.8
$ws.Range("J86").Value = 1350
$ws.Range("K86").Value = 1990.8
$ws.Range("L86").Value = 1350
$ws.Range("M86").Value = -867.8
$ws.Range("N86").Value = -3596

$ws.Range("H89").Value = 1835.4546
$ws.Range("I89").Value = 1990.8
$ws.Range("J89").Value = 1350
$ws.Range("K89").Value = 9954
$ws.Range("L89").Value = 6750
$ws.Range("M89").Value = -4338
$ws.Range("N89").Value = -17982

$ws.Range("H107").Value = 4533.65
$ws.Range("I107").Value = 5510.25
$ws.Range("J107").Value = 3068.75
$ws.Range("K107").Value = 5510.25
$ws.Range("L107").Value = 3068.75
$ws.Range("M107").Value = -3590.25
$ws.Range("N107").Value = -6908.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1675.75
$ws.Range("I31").Value = 1404.9788
$ws.Range("K31").Value = 1404.9788
$ws.Range("M31").Value = -1109.9788

$ws.Range("H34").Value = 1675.75
$ws.Range("I34").Value = 1404.9788
$ws.Range("K34").Value = 1404.9788
$ws.Range("M34").Value = -1202.9788

$ws.Range("H62").Value = 2383.9473
$ws.Range("I62").Value = 2364.4119
$ws.Range("J62").Value = 2550
$ws.Range("K62").Value = 2364.4119
$ws.Range("L62").Value = 2550
$ws.Range("M62").Value = -1740.4119
$ws.Range("N62").Value = -3798

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H65").Value = 2383.9473
$ws.Range("I65").Value = 2364.4119
$ws.Range("J65").Value = 2550
$ws.Range("K65").Value = 11822.0595
$ws.Range("L65").Value = 12750
$ws.Range("M65").Value = -8702.059499999999
$ws.Range("N65").Value = -18990

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100029.4
$ws.Range("I2").Value = 250017.75
$ws.Range("J2").Value = 37.166668
$ws.Range("K2").Value = 1500106.5
$ws.Range("L2").Value = 223.000008
$ws.Range("M2").Value = -1499993.5
$ws.Range("N2").Value = -449.000008

$ws.Range("H95").Value = 6130.769
$ws.Range("J95").Value = 6333.3335
$ws.Range("L95").Value = 19000.0005
$ws.Range("N95").Value = -23118.0005

$ws.Range("H137").Value = 2399.0908
$ws.Range("I137").Value = 1987.7778
$ws.Range("J137").Value = 4250
$ws.Range("K137").Value = 5963.3334
$ws.Range("L137").Value = 12750
$ws.Range("M137").Value = -863.3334000000004
$ws.Range("N137").Value = -22950

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 24114.648
$ws.Range("I70").Value = 32856.51
$ws.Range("J70").Value = 5088.2354
$ws.Range("K70").Value = 32856.51
$ws.Range("L70").Value = 5088.2354
$ws.Range("M70").Value = -32586.51
$ws.Range("N70").Value = -5628.2354

$ws.Range("H73").Value = 24114.648
$ws.Range("I73").Value = 32856.51
$ws.Range("J73").Value = 5088.2354
$ws.Range("K73").Value = 32856.51
$ws.Range("L73").Value = 5088.2354
$ws.Range("M73").Value = -31920.51
$ws.Range("N73").Value = -6960.2354

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1336.8636
$ws.Range("I46").Value = 1146.0769
$ws.Range("J46").Value = 1612.4445
$ws.Range("K46").Value = 1146.0769
$ws.Range("L46").Value = 1612.4445
$ws.Range("M46").Value = -958.0769
$ws.Range("N46").Value = -1988.4445

$ws.Range("H82").Value = 3334
$ws.Range("I82").Value = 2501
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 2501
$ws.Range("L82").Value = 5000
$ws.Range("M82").Value = -2140
$ws.Range("N82").Value = -5722

$ws.Range("H85").Value = 3334
$ws.Range("I85").Value = 2501
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 2501
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = -1253
$ws.Range("N85").Value = -7496

$ws.Range("H122").Value = 3137.348
$ws.Range("I122").Value = 2523
$ws.Range("J122").Value = 3936
$ws.Range("K122").Value = 7569
$ws.Range("L122").Value = 11808
$ws.Range("M122").Value = -5119
$ws.Range("N122").Value = -16708

$ws.Range("H132").Value = 2500.6287
$ws.Range("I132").Value = 2166.25
$ws.Range("J132").Value = 2946.4666
$ws.Range("K132").Value = 6498.75
$ws.Range("L132").Value = 8839.399800000001
$ws.Range("M132").Value = -3968.75
$ws.Range("N132").Value = -13899.3998

$ws.Range("H133").Value = 46342.832
$ws.Range("J133").Value = 46342.832
$ws.Range("L133").Value = 46342.832
$ws.Range("N133").Value = -51402.832

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 7333.3335
$ws.Range("I51").Value = 7000
$ws.Range("K51").Value = 7000
$ws.Range("M51").Value = -6490

$ws.Range("H52").Value = 17995
$ws.Range("J52").Value = 17995
$ws.Range("L52").Value = 17995
$ws.Range("N52").Value = -18447

$ws.Range("H58").Value = 8750
$ws.Range("I58").Value = 9500
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 9500
$ws.Range("L58").Value = 8000
$ws.Range("M58").Value = -9192
$ws.Range("N58").Value = -8616
